$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7: JST SH 4 position 1 mm pitch connector
$ws.Range("A7").Value = "JST SH 4 position 1 mm pitch connector"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 3.95
$ws.Range("D7").Formula = "=B7*C7"
$ws.Range("E7").Value = "1528-4208-ND"
$ws.Range("F7").Value = "This might be a 10 pack (listed on Adafruit, but not sure)"
$ws.Range("G7").Value = "https://www.digikey.com/en/products/detail/adafruit-industries-llc/4208/10230005"
$ws.Range("H7").Value = "https://cdn-shop.adafruit.com/product-files/4208/4208_Kaweei_C13396_diagram.pdf"

# Leave the cursor where the author left it after the edit
$ws.Range("D16").Select()
